$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Mary Lou"
$ws.Range("B4").Value = "Ron"
$ws.Range("C4").Value = "Josh"

$ws.Range("A5").Value = "Hannah"
$ws.Range("B5").Value = "Josh"

$ws.Range("A6").Value = "Daniel"
$ws.Range("B6").Value = "Aaron"

$ws.Range("A7").Select()
